$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "UI for Create Staff TimeTable, Appointment"
$ws.Range("B3").Value = "UI for Login, Avilibilty TimeTable"

$ws.Range("B3").Select() | Out-Null
